$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 (shifts existing rows 18-31 down to 19-32)
$ws.Rows.Item(18).Insert()

# Copy row 17 (unshifted, still intact above the insertion point) into the
# freshly inserted row 18 so formatting and the static/common columns
# (A, B, C, E, F, G, H, I, N, O, Q, R) are duplicated automatically.
$ws.Rows.Item(17).Copy()
$ws.Rows.Item(18).PasteSpecial()

# Overwrite the values that differ for the new row 18, per the diff.
$ws.Cells.Item(18, 4).Value = 44705    # D18  Fecha
$ws.Cells.Item(18, 10).Value = 35      # J18  Volumen
$ws.Cells.Item(18, 11).Value = 26000   # K18  Precio minimo
$ws.Cells.Item(18, 12).Value = 26000   # L18  Precio maximo
$ws.Cells.Item(18, 13).Value = 26000   # M18  Precio promedio ponderado
$ws.Cells.Item(18, 16).Value = 1733    # P18  Precio $/Kg

Write-Host "Inserted new row 18 and shifted subsequent rows down to row 32."
